$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Edit 1: rewrite paragraph beginning "In deriving an optimal spatial
# arrangement..." with the expanded / reworded text.
# ------------------------------------------------------------------
$p20 = $d.Paragraphs.Item(20)
if (-not $p20.Range.Text.StartsWith("In deriving an optimal spatial arrangement")) {
    throw "Paragraph 20 text mismatch: " + $p20.Range.Text.Substring(0, 60)
}
$p20start = $p20.Range.Start
$p20end = $p20.Range.End
$r20 = $d.Range($p20start, $p20end - 1)
$r20.Text = "In deriving an optimal spatial arrangement of nodes for this scenario, an intelligent strategy needs to be implemented to maximize the total volume enclosed, while eliminating or at least minimizing shadow zones. We choose a genetic algorithm-based approach to solving this problem. Genetic algorithms provide an evolutionary approach towards solving such problems by aiming to improve the fitness of each successive generation, mimicking the evolution of life in nature on a suitably simplistic scale. An initial population of individuals is required, and is often randomly seeded. A fitness function is defined, which assigns a score to every member of the current population based on the evaluation of relevant characteristics. The fittest individuals from this pool are selected for breeding to create the next generation. Additional factors like mutation in chromosomes and crossovers among subpopulations may also be specified to reduce the chances of the solution converging towards a local maximum. Since shadow zones in such a scenario are essentially holes in the coverage shell, penalties are required to discourage such arrangements from participating in the evolution of the genome. In every successive generation, the score of the best-fit individual is expected to improve due to selective breeding. As the score stagnates with respect to average change in fitness, generation, or time, the algorithm terminates with the optimal solution as its output."

# ------------------------------------------------------------------
# Edit 2: merge the next two paragraphs ("In our deployment
# strategy..." and "For every individual, a convex hull...") into a
# single paragraph carrying the new, reworded / reorganized text.
# ------------------------------------------------------------------
$p21 = $d.Paragraphs.Item(21)
$p22 = $d.Paragraphs.Item(22)
if (-not $p21.Range.Text.StartsWith("In our deployment strategy")) {
    throw "Paragraph 21 text mismatch: " + $p21.Range.Text.Substring(0, 60)
}
if (-not $p22.Range.Text.StartsWith("For every individual, a convex hull")) {
    throw "Paragraph 22 text mismatch: " + $p22.Range.Text.Substring(0, 60)
}
$p21start = $p21.Range.Start
$p22end = $p22.Range.End

# Remove the paragraph mark that separates the two paragraphs so they
# become a single paragraph (keeping paragraph 21's pPr/formatting).
$markRange = $d.Range($p21.Range.End - 1, $p21.Range.End)
$markRange.Delete()

# The merged paragraph now spans p21start .. (p22end - 1); p22end - 1
# stops just before the (still present) trailing paragraph mark of the
# merged paragraph.
$mergedEnd = $p22end - 1
$r21 = $d.Range($p21start, $mergedEnd)
$r21.Text = "In our deployment strategy, we model a given set of sensor nodes as point sources with specified characteristics, and allow them to spread the spatial arrangement using a genetic algorithm. We create the initial population by randomly scattering the nodes in three dimensions such that every node is initially within range of every other node. In a real life situation, given the upper and lower depth bounds, the goal would be to secure the disk-like chunk in the slab of the water body. With each successive generation, the node arrangement is allowed to expand in the desired aspect ratio. For every individual, a convex hull is stretched over the point cloud formed by the nodes in three dimensional space to form a polyhedron. The volume of this polyhedron not only serves as the initial score for the individual prior to constraint checking, but its visualization can also be used to identify shadow zones as well as highlight nodes surplus to requirements in achieving the given objective. Shadow zones can be born of holes in the facets of the polyhedron, or gaps in the edge coverage, which are a lot worse. To eliminate these gaps, every edge of every facet of the polyhedron is tested for overlap-accounted total edge coverage. As soon as the first violation of this constraint is found, the individual is rejected outright with a score of zero. Surplus nodes lie inside the polyhedron, and at best might possibly contribute to face coverage and alternate communication routes. In the case of the ultimate best-fit individual, a cost-benefit analysis of the solution might serve to decide the involvement of these nodes in the actual deployment."

# ------------------------------------------------------------------
# Edit 3: the paragraph "Since shadow zones in such a scenario..."
# (now the final paragraph in the document, right after the
# "echo-based detection range..." paragraph) is removed outright -
# its content was folded into the rewritten paragraph above.
# ------------------------------------------------------------------
$lastIndex = $d.Paragraphs.Count
$pLast = $d.Paragraphs.Item($lastIndex)
if (-not $pLast.Range.Text.StartsWith("Since shadow zones in such a scenario")) {
    throw "Last paragraph text mismatch: " + $pLast.Range.Text.Substring(0, 60)
}
$delRange = $d.Range($pLast.Range.Start, $pLast.Range.End)
$delRange.Delete()

Write-Output ("Edits applied. Paragraph count now: " + $d.Paragraphs.Count)
